# example 005 - add STAT.STAT_ID / STAT.VALUE columns (G, H) to the Item sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows -----------------------------------------------------
# Row 1: field names, Row 2: field types (matches the existing A:F header pattern)
$ws.Range("G1").Value = "STAT.STAT_ID"
$ws.Range("H1").Value = "STAT.VALUE"
$ws.Range("G2").Value = "INT"
$ws.Range("H2").Value = "FLOAT"

# --- Data rows ---------------------------------------------------------
# Rows 3-32 (the 30 item records): new columns default to 0.
for ($r = 3; $r -le 32; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = 0
}

# --- Column sizing -------------------------------------------------------
$ws.Columns("G").ColumnWidth = 13.375

# --- Selection -----------------------------------------------------------
[void]$ws.Range("H3").Select()
